$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new values look like plain numbers need to be pinned
# to Text format first, otherwise Excel auto-converts them to numeric cells
# (e.g. "0.0000241" -> 2.41E-05), which does not match the source data which
# always stores these coin prices as text.
$numericLookingDCells = @(5,6,10,14,15,16,20,21,22,24,25,26,27,28,33,34,35,37,40,41,42,43,44,46,47,49,50) | ForEach-Object { "D$_" }
foreach ($addr in $numericLookingDCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "94.177.73"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").Value = "3.071.84"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "235.28"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "610.58"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "0.803"
$ws.Range("E10").Value = "  +8.21%  "

$ws.Range("D11").Value = "3.070.10"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("D13").Value = "93.975.58"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").Value = "  -2.47%  "

$ws.Range("D15").Value = "33.64"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  -1.81%  "

$ws.Range("D17").Value = "3.646.95"
$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").Value = "3.055.29"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("E19").Value = "  -4.97%  "

$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  -1.87%  "

$ws.Range("D21").Value = "5.66"
$ws.Range("E21").Value = "  -0.99%  "

$ws.Range("D22").Value = "438.13"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("E23").Value = "  -4.43%  "

$ws.Range("D24").Value = "0.0000190"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").Value = "8.34"
$ws.Range("E25").Value = "  +6.30%  "

$ws.Range("D26").Value = "5.52"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("D27").Value = "84.75"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").Value = "11.85"
$ws.Range("E28").Value = "  +2.28%  "

$ws.Range("D29").Value = "3.229.20"
$ws.Range("E29").Value = "  -1.22%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("E31").Value = "  +8.34%  "

$ws.Range("E32").Value = "  +6.45%  "

$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -7.68%  "

$ws.Range("D34").Value = "9.03"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "7.78"
$ws.Range("E35").Value = "  -1.59%  "

$ws.Range("E36").Value = "  -2.62%  "

$ws.Range("D37").Value = "25.39"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("E38").Value = "  -11.24%  "

$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("D40").Value = "24.01"
$ws.Range("E40").Value = "  +0.65%  "

$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").Value = "3.76"
$ws.Range("E41").Value = "  -3.21%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "472.44"
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.435"
$ws.Range("E43").Value = "  +1.68%  "

$ws.Range("D44").Value = "1.27"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").Value = "3.13"
$ws.Range("E46").Value = "  -4.75%  "

$ws.Range("D47").Value = "161.59"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("E48").Value = "  -1.47%  "

$ws.Range("D49").Value = "1.82"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").Value = "43.63"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("E51").Value = "  +0.03%  "
